# Auto-generated edit script applying numeric corrections to Titan_Profits data
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 211.92308
$ws.Range("I5").Value = 31.625
$ws.Range("K5").Value = 31.625
$ws.Range("M5").Value = 83.375

$ws.Range("H10").Value = 7500
$ws.Range("J10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("N10").Value = -20586

$ws.Range("H39").Value = 235.05882
$ws.Range("I39").Value = 76.666664
$ws.Range("J39").Value = 615.2
$ws.Range("K39").Value = 229.999992
$ws.Range("L39").Value = 1845.6
$ws.Range("M39").Value = 66.00000800000001
$ws.Range("N39").Value = -2437.6

$ws.Range("H43").Value = 1201
$ws.Range("I43").Value = 1201
$ws.Range("K43").Value = 1201
$ws.Range("M43").Value = -1132

$ws.Range("H64").Value = 838581.8
$ws.Range("I64").Value = 1114664.6
$ws.Range("J64").Value = 10333.333
$ws.Range("K64").Value = 1114664.6
$ws.Range("L64").Value = 10333.333
$ws.Range("M64").Value = -1114416.6
$ws.Range("N64").Value = -10829.333

$ws.Range("H67").Value = 838581.8
$ws.Range("I67").Value = 1114664.6
$ws.Range("J67").Value = 10333.333
$ws.Range("K67").Value = 1114664.6
$ws.Range("L67").Value = 10333.333
$ws.Range("M67").Value = -1113806.6
$ws.Range("N67").Value = -12049.333

$ws.Range("H76").Value = 5850887.5
$ws.Range("I76").Value = 6538768.5
$ws.Range("J76").Value = 3900
$ws.Range("K76").Value = 6538768.5
$ws.Range("L76").Value = 3900
$ws.Range("M76").Value = -6538453.5
$ws.Range("N76").Value = -4530

$ws.Range("H79").Value = 5850887.5
$ws.Range("I79").Value = 6538768.5
$ws.Range("J79").Value = 3900
$ws.Range("K79").Value = 6538768.5
$ws.Range("L79").Value = 3900
$ws.Range("M79").Value = -6537676.5
$ws.Range("N79").Value = -6084

$ws.Range("H113").Value = 57505.58
$ws.Range("I113").Value = 72180.664
$ws.Range("J113").Value = 2474
$ws.Range("K113").Value = 72180.664
$ws.Range("L113").Value = 2474
$ws.Range("M113").Value = -68926.664
$ws.Range("N113").Value = -8982

$ws.Range("H132").Value = 330305.94
$ws.Range("I132").Value = 405914.2
$ws.Range("K132").Value = 1217742.6
$ws.Range("M132").Value = -1215212.6

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17336.463
$ws.Range("I32").Value = 2938.15
$ws.Range("K32").Value = 2938.15
$ws.Range("M32").Value = -2651.15

$ws.Range("H61").Value = 2642.05
$ws.Range("I61").Value = 2127.5
$ws.Range("K61").Value = 2127.5
$ws.Range("M61").Value = -1915.5

$ws.Range("H132").Value = 2460.1177
$ws.Range("I132").Value = 2170.238
$ws.Range("K132").Value = 6510.714
$ws.Range("M132").Value = -3980.714

$ws.Range("H136").Value = 2642.05
$ws.Range("I136").Value = 2127.5
$ws.Range("K136").Value = 6382.5
$ws.Range("M136").Value = -3832.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 361
$ws.Range("I11").Value = 400
$ws.Range("J11").Value = 302.5
$ws.Range("K11").Value = 400
$ws.Range("L11").Value = 302.5
$ws.Range("M11").Value = -260
$ws.Range("N11").Value = -582.5

$ws.Range("H134").Value = 5119.85
$ws.Range("I134").Value = 4439.3335
$ws.Range("J134").Value = 5411.5
$ws.Range("K134").Value = 13318.0005
$ws.Range("L134").Value = 16234.5
$ws.Range("M134").Value = -10783.0005
$ws.Range("N134").Value = -21304.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1291.2424
$ws.Range("I31").Value = 879
$ws.Range("K31").Value = 879
$ws.Range("M31").Value = -584

$ws.Range("H32").Value = 15166.333
$ws.Range("I32").Value = 15166.333
$ws.Range("K32").Value = 15166.333
$ws.Range("M32").Value = -14850.333

$ws.Range("H34").Value = 1291.2424
$ws.Range("I34").Value = 879
$ws.Range("K34").Value = 879
$ws.Range("M34").Value = -677

$ws.Range("H132").Value = 2236.3125
$ws.Range("I132").Value = 1898.7222
$ws.Range("J132").Value = 3249.0833
$ws.Range("K132").Value = 5696.1666
$ws.Range("L132").Value = 9747.249899999999
$ws.Range("M132").Value = -3166.1666
$ws.Range("N132").Value = -14807.2499

$ws.Range("H134").Value = 3304.2144
$ws.Range("I134").Value = 1807.7333
$ws.Range("K134").Value = 5423.199900000001
$ws.Range("M134").Value = -2888.199900000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1214.1428
$ws.Range("I116").Value = 999.8333
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 2999.4999
$ws.Range("L116").Value = 7500
$ws.Range("M116").Value = 442.5001000000002
$ws.Range("N116").Value = -14384

$ws.Range("H139").Value = 1817
$ws.Range("I139").Value = 1506.6
$ws.Range("K139").Value = 4519.799999999999
$ws.Range("M139").Value = 620.2000000000007

$ws.Range("H141").Value = 4390.9
$ws.Range("I141").Value = 5201.2856
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 15603.8568
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = -10423.8568
$ws.Range("N141").Value = -17860

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 6890
$ws.Range("J17").Value = 2102
$ws.Range("L17").Value = 2102
$ws.Range("N17").Value = -2438

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H80").Value = 2472
$ws.Range("I80").Value = 2434.889
$ws.Range("J80").Value = 2513.75
$ws.Range("K80").Value = 2434.889
$ws.Range("L80").Value = 2513.75
$ws.Range("M80").Value = -1436.889
$ws.Range("N80").Value = -4509.75

$ws.Range("H83").Value = 2472
$ws.Range("I83").Value = 2434.889
$ws.Range("J83").Value = 2513.75
$ws.Range("K83").Value = 12174.445
$ws.Range("L83").Value = 12568.75
$ws.Range("M83").Value = -7182.445
$ws.Range("N83").Value = -22552.75

$ws.Range("H122").Value = 1853736.5
$ws.Range("I122").Value = 3705003.8
$ws.Range("K122").Value = 11115011.4
$ws.Range("M122").Value = -11112561.4

$ws.Range("H132").Value = 3110.84
$ws.Range("I132").Value = 2982.3225
$ws.Range("J132").Value = 3320.5264
$ws.Range("K132").Value = 8946.967500000001
$ws.Range("L132").Value = 9961.5792
$ws.Range("M132").Value = -6416.967500000001
$ws.Range("N132").Value = -15021.5792

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 879.41174
$ws.Range("I22").Value = 537.5
$ws.Range("J22").Value = 1183.3334
$ws.Range("K22").Value = 537.5
$ws.Range("L22").Value = 1183.3334
$ws.Range("M22").Value = -242.5
$ws.Range("N22").Value = -1773.3334

$ws.Range("H27").Value = 879.41174
$ws.Range("I27").Value = 537.5
$ws.Range("J27").Value = 1183.3334
$ws.Range("K27").Value = 537.5
$ws.Range("L27").Value = 1183.3334
$ws.Range("M27").Value = -430.5
$ws.Range("N27").Value = -1397.3334

$ws.Range("H46").Value = 1963.5555
$ws.Range("I46").Value = 1042.5
$ws.Range("J46").Value = 2700.4
$ws.Range("K46").Value = 1042.5
$ws.Range("L46").Value = 2700.4
$ws.Range("M46").Value = -854.5
$ws.Range("N46").Value = -3076.4

$ws.Range("H55").Value = 353.2857
$ws.Range("I55").Value = 273.875
$ws.Range("J55").Value = 459.16666
$ws.Range("K55").Value = 273.875
$ws.Range("L55").Value = 459.16666
$ws.Range("M55").Value = -100.875
$ws.Range("N55").Value = -805.16666

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H107").Value = 891.30304
$ws.Range("I107").Value = 746.5217
$ws.Range("J107").Value = 1224.3
$ws.Range("K107").Value = 2239.5651
$ws.Range("L107").Value = 3672.9
$ws.Range("M107").Value = -319.5650999999998
$ws.Range("N107").Value = -7512.9

$ws.Range("H126").Value = 92186.09
$ws.Range("I126").Value = 111960.78
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 335882.34
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -333412.34
$ws.Range("N126").Value = -14540

$ws.Range("H132").Value = 16131211
$ws.Range("I132").Value = 19232642
$ws.Range("K132").Value = 57697926
$ws.Range("M132").Value = -57695396
